# Trade #22 closed at 2026-02-17 20:07:05 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.81   # Current Capital
$summary.Range("B4").Value = -0.2      # Total P&L $
$summary.Range("B5").Value = -0.18     # Total P&L %
$summary.Range("B6").Value = 22        # Total Trades
$summary.Range("B7").Value = 11        # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.81      # Capital
$status.Range("D5").Value = 22         # Trades
$status.Range("E5").Value = -0.2       # P&L $
$status.Range("F5").Value = -0.19      # P&L %
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------
# New trade row (#22 / spreadsheet row 23) appended to both the
# "All Trades" sheet and the "MarketMaking" sheet.
# ---------------------------------------------------------------
function Add-TradeRow23($ws) {
    $ws.Range("A23").Value = 22

    # Date / time columns are plain text in this workbook, not real
    # Excel dates - force text number format first so COM does not
    # silently convert the "2026-02-17" string into a date serial.
    $ws.Range("B23").NumberFormat = "@"
    $ws.Range("B23").Value = "2026-02-17"
    $ws.Range("B23").ClearFormats()

    $ws.Range("C23").Value = "20:06:59"
    $ws.Range("D23").Value = "MarketMaking"
    $ws.Range("E23").Value = "UP"
    $ws.Range("F23").Value = 0.77
    $ws.Range("G23").Value = 0.78
    $ws.Range("H23").Value = "CLOSED"
    $ws.Range("I23").Value = 1.2987
    $ws.Range("J23").Value = 0.01
    $ws.Range("K23").Value = 99.81
    $ws.Range("L23").Value = 0
    $ws.Range("M23").Value = 0
    $ws.Range("N23").Value = 0.6
    $ws.Range("O23").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P23").Value = "early_exit"
    $ws.Range("Q23").Value = 0.12
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow23 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow23 $marketMaking
